# Parametrização de pesos de carga, spread e pico máximo de turmas permitido
# Updates the distribution counts in the "Sheet1" consolidation table
# (columns B..J, rows 2..19) to reflect the re-weighted allocation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# PROG_1 (row 2)
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 0

# PROG_10 (row 3)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

# PROG_11 (row 4)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 3

# PROG_12 (row 5)
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 14

# PROG_13 (row 6)
$ws.Range("D6").Value = 4
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 0

# PROG_2 (row 7)
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 4

# PROG_3 (row 8)
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 4
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 14

# PROG_4 (row 9)
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 13

# PROG_5 (row 10)
$ws.Range("C10").Value = 1
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 4

# PROG_6 (row 11)
$ws.Range("D11").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 13

# PROG_7 (row 12)
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 3
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 3

# PROG_8 (row 13)
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 14

# PROG_9 (row 14)
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 13

# ROB_1 (row 15)
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 12

# ROB_2 (row 16)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 2
$ws.Range("J16").Value = 11

# ROB_3 (row 17)
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 11

# ROB_4 (row 18)
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 12

# ROB_5 (row 19)
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 2
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 5
